$d = $word.ActiveDocument

# Locate the bibliography paragraph (the one that starts with "WELLER, Mark").
$targetParagraph = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("WELLER, Mark")) {
        $targetParagraph = $p
        break
    }
}

if ($targetParagraph -eq $null) {
    Write-Host "ERROR: bibliography paragraph (starting with 'WELLER, Mark') was not found."
}

# The bibliography used to be one giant run of concatenated references.
# Turn it into the same text broken up by a blank line (two manual line
# breaks, ^l^l) between each reference. Each boundary string is unique in
# the document, so a literal (non-wildcard) Find/Replace anchored on the
# text spanning the seam reliably inserts the breaks in the right spot.
# Find is re-run against a freshly fetched (and now longer) paragraph range
# each time, with Wrap forced to wdFindStop (0) so it can never search
# outside this paragraph.

$replacements = @(
    @{ Find = "2017. E-book. CHANG, Raymond"; Replace = "2017. E-book. ^l^lCHANG, Raymond" },
    @{ Find = "AMGH Editora Ltda., 2010.BROWN, T.L."; Replace = "AMGH Editora Ltda., 2010.^l^lBROWN, T.L." },
    @{ Find = "Pearson Prentice Hall, 2005-2007.BRADY, J"; Replace = "Pearson Prentice Hall, 2005-2007.^l^lBRADY, J" },
    @{ Find = "Técnicos Científicos, 1981.LEE, J. D.,"; Replace = "Técnicos Científicos, 1981.^l^lLEE, J. D.," },
    @{ Find = "Ltda. SP-2001.SHRIVER, D."; Replace = "Ltda. SP-2001.^l^lSHRIVER, D." },
    @{ Find = "Porto Alegre-RS, 2008.QUAGLIANO, J.V"; Replace = "Porto Alegre-RS, 2008.^l^lQUAGLIANO, J.V" }
)

foreach ($item in $replacements) {
    $rng = $targetParagraph.Range
    $found = $rng.Find.Execute($item.Find, $false, $false, $false, $false, $false, $true, 0, $false, $item.Replace, 2)
    if (-not $found) {
        Write-Host "WARNING: could not find/replace boundary text: " $item.Find
    }
}
